$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The data table (Periodo Mora / Valor Mora) lives in rows 16-50, columns E (5) and F (6).
# This edit reverses the row order of that period/value list: the most recent period
# (1912) now appears first (row 16) together with its original partial value (29867),
# while the oldest period (1702) now appears last (row 50) with the standard value (56000).

$firstRow = 16
$lastRow = 50

# Snapshot the current Periodo Mora (col E) and Valor Mora (col F) values before
# overwriting anything, so the reversal is computed from a stable source.
$periods = @()
$values = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
    $values += $ws.Cells.Item($r, 6).Value2
}

$count = $lastRow - $firstRow + 1
for ($i = 0; $i -lt $count; $i++) {
    $targetRow = $firstRow + $i
    $sourceIndex = $count - 1 - $i
    $ws.Cells.Item($targetRow, 5).Value = $periods[$sourceIndex]
    $ws.Cells.Item($targetRow, 6).Value = $values[$sourceIndex]
}
